$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# Update price (D) and volume-change (E) columns for rows 2-49
Set-TextValue $ws.Range("D2") "25.890.88"
$ws.Range("E2").Value = "  -0.05%  "
Set-TextValue $ws.Range("D3") "1.731.52"
$ws.Range("E3").Value = "  -0.62%  "
Set-TextValue $ws.Range("D4") "1.001"
$ws.Range("E4").Value = "  +0.13%  "
Set-TextValue $ws.Range("D5") "245.29"
$ws.Range("E5").Value = "  +2.69%  "
Set-TextValue $ws.Range("D6") "1.001"
$ws.Range("E6").Value = "  +0.11%  "
Set-TextValue $ws.Range("D7") "0.5018"
$ws.Range("E7").Value = "  -3.13%  "
Set-TextValue $ws.Range("D8") "0.2723"
$ws.Range("E8").Value = "  -0.86%  "
Set-TextValue $ws.Range("D9") "0.06166"
$ws.Range("E9").Value = "  +0.33%  "
Set-TextValue $ws.Range("D10") "1.741.01"
$ws.Range("E10").Value = "  -0.05%  "
Set-TextValue $ws.Range("D11") "0.07236"
$ws.Range("E11").Value = "  +0.90%  "
Set-TextValue $ws.Range("D12") "0.6517"
$ws.Range("E12").Value = "  +1.06%  "
Set-TextValue $ws.Range("D13") "15.14"
$ws.Range("E13").Value = "  +1.01%  "
Set-TextValue $ws.Range("D14") "4.770"
$ws.Range("E14").Value = "  +3.74%  "
Set-TextValue $ws.Range("D15") "76.99"
$ws.Range("E15").Value = "  -0.62%  "
Set-TextValue $ws.Range("D16") "0.9987"
$ws.Range("E16").Value = "  -0.12%  "
Set-TextValue $ws.Range("D17") "1.001"
$ws.Range("E17").Value = "  +0.20%  "
Set-TextValue $ws.Range("D18") "25.897.69"
$ws.Range("E18").Value = "  -0.06%  "
Set-TextValue $ws.Range("D19") "11.91"
$ws.Range("E19").Value = "  +1.60%  "
Set-TextValue $ws.Range("D20") "0.000006817"
$ws.Range("E20").Value = "  +0.50%  "
Set-TextValue $ws.Range("D21") "4.594"
$ws.Range("E21").Value = "  +7.41%  "
Set-TextValue $ws.Range("D22") "1.962.63"
$ws.Range("E22").Value = "  +0.01%  "
Set-TextValue $ws.Range("D23") "8.803"
$ws.Range("E23").Value = "  +1.49%  "
Set-TextValue $ws.Range("D24") "5.479"
$ws.Range("E24").Value = "  +4.43%  "
Set-TextValue $ws.Range("D25") "133.91"
$ws.Range("E25").Value = "  -3.47%  "
Set-TextValue $ws.Range("D26") "15.28"
$ws.Range("E26").Value = "  +0.96%  "
Set-TextValue $ws.Range("D27") "1.790"
$ws.Range("E27").Value = "  +1.42%  "
Set-TextValue $ws.Range("D28") "1.421"
$ws.Range("E28").Value = "  -5.91%  "
Set-TextValue $ws.Range("D29") "105.62"
$ws.Range("E29").Value = "  -0.39%  "
Set-TextValue $ws.Range("D30") "3.977"
$ws.Range("E30").Value = "  +0.96%  "
Set-TextValue $ws.Range("D31") "0.08130"
$ws.Range("E31").Value = "  -2.11%  "
Set-TextValue $ws.Range("D32") "3.687"
$ws.Range("E32").Value = "  +0.97%  "
Set-TextValue $ws.Range("D33") "0.04729"
$ws.Range("E33").Value = "  +2.97%  "
Set-TextValue $ws.Range("D34") "2.652"
$ws.Range("E34").Value = "  -0.18%  "
Set-TextValue $ws.Range("D35") "0.9954"
$ws.Range("E35").Value = "  +0.47%  "
Set-TextValue $ws.Range("D36") "0.6133"
$ws.Range("E36").Value = "  -0.87%  "
Set-TextValue $ws.Range("D37") "2.742"
$ws.Range("E37").Value = "  +1.97%  "
Set-TextValue $ws.Range("D38") "0.01608"
$ws.Range("E38").Value = "  -0.30%  "
Set-TextValue $ws.Range("D39") "0.8805"
$ws.Range("E39").Value = "  +18.98%  "
Set-TextValue $ws.Range("D40") "1.955"
$ws.Range("E40").Value = "  +1.07%  "
Set-TextValue $ws.Range("D41") "1.000"
$ws.Range("E41").Value = "  +0.08%  "
Set-TextValue $ws.Range("D42") "101.40"
$ws.Range("E42").Value = "  +3.51%  "
Set-TextValue $ws.Range("D43") "0.3901"
$ws.Range("E43").Value = "  +1.61%  "
Set-TextValue $ws.Range("D44") "5.009"
$ws.Range("E44").Value = "  +0.54%  "
Set-TextValue $ws.Range("D45") "0.1178"
$ws.Range("E45").Value = "  +4.55%  "
Set-TextValue $ws.Range("D46") "6.358"
$ws.Range("E46").Value = "  +2.32%  "
Set-TextValue $ws.Range("D47") "55.68"
$ws.Range("E47").Value = "  +1.50%  "
Set-TextValue $ws.Range("D48") "0.05280"
$ws.Range("E48").Value = "  +0.40%  "
Set-TextValue $ws.Range("D49") "30.77"
$ws.Range("E49").Value = "  +0.78%  "

# Row 50/51: EnergySwap and Decentraland swap positions, with updated price/volume data
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws.Range("D50") "7.654"
$ws.Range("E50").Value = "  +0.73%  "

$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue $ws.Range("D51") "0.3485"
$ws.Range("E51").Value = "  +2.31%  "
